# Auto-generated Excel COM-interop script
# Applies the "Update countries & provincias Spain" edit to the Pais sheet:
#   - Refreshed case counts (as of 1 Mayo 2020, 09:22) for several country rows
#   - Singapur/Pakistan swapped ranking (rows 27-28)
#   - Armenia moved up past Irak/Croacia/Ghana (rows 68-71)
#   - Updated "Datos actualizados" timestamp in A1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Country-name / label cells (text) ---
$ws.Range("A1").Value = "Datos actualizados a 1 de Mayo de 2020 a las 09:22"
$ws.Range("A27").Value = "Singapur"
$ws.Range("A28").Value = "Pakistan"
$ws.Range("A68").Value = "Armenia"
$ws.Range("A69").Value = "Irak"
$ws.Range("A70").Value = "Croacia"
$ws.Range("A71").Value = "Ghana"

# --- Updated numeric statistics ---
$ws.Range("B27").Value = 17101
$ws.Range("C27").Value = 932
$ws.Range("D27").Value = 1244
$ws.Range("E27").Value = 15842
$ws.Range("F27").Value = 21
$ws.Range("G27").Value = 0
$ws.Range("H27").Value = 15
$ws.Range("B28").Value = 16817
$ws.Range("C28").Value = 344
$ws.Range("D28").Value = 4315
$ws.Range("E28").Value = 12117
$ws.Range("F28").Value = 111
$ws.Range("G28").Value = 24
$ws.Range("H28").Value = 385
$ws.Range("B30").Value = 16004
$ws.Range("C30").Value = 58
$ws.Range("D30").Value = 8758
$ws.Range("E30").Value = 7023
$ws.Range("G30").Value = 1
$ws.Range("H30").Value = 223
$ws.Range("B45").Value = 7689
$ws.Range("C45").Value = 7
$ws.Range("E45").Value = 4138
$ws.Range("F45").Value = 67
$ws.Range("G45").Value = 1
$ws.Range("H45").Value = 237
$ws.Range("B61").Value = 3551
$ws.Range("C61").Value = 149
$ws.Range("E61").Value = 2660
$ws.Range("F61").Value = 40
$ws.Range("B68").Value = 2148
$ws.Range("C68").Value = 82
$ws.Range("D68").Value = 977
$ws.Range("E68").Value = 1138
$ws.Range("F68").Value = 10
$ws.Range("G68").Value = 1
$ws.Range("H68").Value = 33
$ws.Range("B69").Value = 2085
$ws.Range("D69").Value = 1375
$ws.Range("E69").Value = 617
$ws.Range("F69").Value = 0
$ws.Range("H69").Value = 93
$ws.Range("B70").Value = 2076
$ws.Range("D70").Value = 1348
$ws.Range("E70").Value = 659
$ws.Range("F70").Value = 20
$ws.Range("H70").Value = 69
$ws.Range("B71").Value = 2074
$ws.Range("D71").Value = 212
$ws.Range("E71").Value = 1845
$ws.Range("F71").Value = 4
$ws.Range("H71").Value = 17
$ws.Range("B93").Value = 870
$ws.Range("C93").Value = 12
$ws.Range("E93").Value = 506
$ws.Range("F93").Value = 2
$ws.Range("G93").Value = 1
$ws.Range("H93").Value = 16
$ws.Range("D126").Value = 219
$ws.Range("E126").Value = 96
